# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# cells to reflect the newly generated handback report timestamps.
$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 21:04:17"

# "zh-cn" sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 21:04:13"
$wsZhCn.Range("K2").Value = "2016-08-28 21:04:30"

# "de-de" sheet: Correspond Handoff Datetime (H2, shares text with Overview!G2)
# and Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 21:04:17"
$wsDeDe.Range("K2").Value = "2016-08-28 21:04:37"
